# Update "想去人数" (interested-count) figures in the "展览" and "全部类型"
# sheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F5"  = 3078
    "F7"  = 2431
    "F9"  = 120
    "F11" = 1226
    "F15" = 1094
    "F16" = 300
    "F21" = 63
    "F23" = 131
    "F25" = 246
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
